$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some cells in columns L/M use a Text ("@") number format, so a
# plain .Value assignment of a number gets coerced to a text string.
# Temporarily switch to a numeric format, assign, then restore "@" so the
# stored cell keeps its original look (Text format) but a real numeric value.
function Set-NumericValueOnTextCell($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "0"
    $rng.Value = $value
    $rng.NumberFormat = "@"
}

# Row 348: new positive cases revised upward
$ws.Range("C348").Value = 47

# Row 353: one extra-hospital death recorded
Set-NumericValueOnTextCell "M353" 1

# Row 355: new positive cases revised upward
$ws.Range("C355").Value = 35

# Row 356: new positive cases revised upward, plus one more hospital death
$ws.Range("C356").Value = 62
Set-NumericValueOnTextCell "L356" 3

# Row 357: new positive cases revised upward, plus one hospital death
$ws.Range("C357").Value = 48
Set-NumericValueOnTextCell "L357" 1

# Row 358: this day's figures were not yet available before, now filled in
$ws.Range("C358").Value = 10
$ws.Range("E358").Value = 9
$ws.Range("F358").Value = 7
$ws.Range("G358").Value = 37
Set-NumericValueOnTextCell "L358" 0
Set-NumericValueOnTextCell "M358" 0

$wb.Save()
